$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value = 1187536.510541698
$ws.Range("G8").Value = 1188408.185213571
$ws.Range("E9").Value = 845075.5262215403
$ws.Range("F9").Value = 529773.4894583018
$ws.Range("G9").Value = 528901.8147864268
